$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("data"): add new column AF ("20. 7. 2021") after AE, matching header/data style ---
$ws1.Range("AE1:AE23").Copy($ws1.Range("AF1:AF23"))
$ws1.Range("AF1").Value = "20. 7. 2021"
$ws1.Range("AF2").Value = 0.1
$ws1.Range("AF3").Value = 0.09
$ws1.Range("AF4").Value = 0.1
$ws1.Range("AF5").Value = 0.06
$ws1.Range("AF6").Value = 0.11
$ws1.Range("AF7").Value = 0.13
$ws1.Range("AF8").Value = 0.09
$ws1.Range("AF9").Value = 0.2
$ws1.Range("AF10").Value = 0.11
$ws1.Range("AF11").Value = 0.08
$ws1.Range("AF12").Value = 0.09
$ws1.Range("AF13").Value = 0.08
$ws1.Range("AF14").Value = 0.19
$ws1.Range("AF15").Value = 0.11
$ws1.Range("AF16").Value = 0.08
$ws1.Range("AF17").Value = 0.15
$ws1.Range("AF18").Value = 0.1
$ws1.Range("AF19").Value = 0.07
$ws1.Range("AF20").Value = 0.09
$ws1.Range("AF21").Value = 0.05
$ws1.Range("AF22").Value = 0.07
$ws1.Range("AF23").Value = 0.18

# --- Sheet2 ("pocetR"): add new column AE ("20. 7. 2021") after AD, matching header/data style ---
$ws2.Range("AD1:AD23").Copy($ws2.Range("AE1:AE23"))
$ws2.Range("AE1").Value = "20. 7. 2021"
$ws2.Range("AE2").Value = 1782
$ws2.Range("AE3").Value = 863
$ws2.Range("AE4").Value = 919
$ws2.Range("AE5").Value = 233
$ws2.Range("AE6").Value = 630
$ws2.Range("AE7").Value = 285
$ws2.Range("AE8").Value = 634
$ws2.Range("AE9").Value = 145
$ws2.Range("AE10").Value = 286
$ws2.Range("AE11").Value = 334
$ws2.Range("AE12").Value = 314
$ws2.Range("AE13").Value = 703
$ws2.Range("AE14").Value = 158
$ws2.Range("AE15").Value = 372
$ws2.Range("AE16").Value = 1252
$ws2.Range("AE17").Value = 201
$ws2.Range("AE18").Value = 684
$ws2.Range("AE19").Value = 557
$ws2.Range("AE20").Value = 236
$ws2.Range("AE21").Value = 459
$ws2.Range("AE22").Value = 767
$ws2.Range("AE23").Value = 556

# --- Update the footer title strings (publish date 28. 6. 2021 -> 27. 7. 2021) ---
$ws1.Range("A24").Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 27. 7. 2021"
$ws2.Range("A24").Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 27. 7. 2021"
